$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to be bumped by
# one day (46074 -> 46075) for every data row (rows 2 through 386).
for ($r = 2; $r -le 386; $r++) {
    $ws.Cells.Item($r, 3).Value = 46075
}
